# Update cryptocurrency price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.589.53"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "'2.357.41"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'331.52"
$ws.Range("E5").Value = "  +6.35%  "
$ws.Range("D6").Value = "'100.17"
$ws.Range("E6").Value = "  -8.37%  "
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.630"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'39.71"
$ws.Range("E10").Value = "  -6.40%  "
$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'8.45"
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("E13").Value = "  -4.39%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'16.32"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "'2.718.26"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'2.363.12"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'42.608.69"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").Value = "'7.93"
$ws.Range("E19").Value = "  +8.62%  "
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("E21").Value = "  +9.45%  "
$ws.Range("D22").Value = "'75.40"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'268.77"
$ws.Range("E23").Value = "  +6.10%  "
$ws.Range("D25").Value = "'9.94"
$ws.Range("E25").Value = "  +10.03%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'11.47"
$ws.Range("E27").Value = "  -4.90%  "
$ws.Range("D28").Value = "'23.13"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").Value = "'176.08"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "'3.08"
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("D32").Value = "'0.0903"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Value = "'35.42"
$ws.Range("E33").Value = "  -9.80%  "
$ws.Range("D34").Value = "'6.11"
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D36").Value = "'4.60"
$ws.Range("E36").Value = "  -8.33%  "
$ws.Range("E37").Value = "  -4.93%  "
$ws.Range("D38").Value = "'2.95"
$ws.Range("E38").Value = "  +9.27%  "
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").Value = "'3.83"
$ws.Range("E40").Value = "  -7.11%  "
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "'70.14"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'117.85"
$ws.Range("E45").Value = "  +6.98%  "
$ws.Range("D46").Value = "'90.75"
$ws.Range("E46").Value = "  +30.37%  "
$ws.Range("D47").Value = "'11.93"
$ws.Range("E47").Value = "  -6.69%  "
$ws.Range("D48").Value = "'5.47"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").Value = "'9.12"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.27"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.572.04"
$ws.Range("E51").Value = "  +5.46%  "
